$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.0157010555267334, 0.001296758651733398, 0.9333333333333333)
    3  = @(0.02471041679382324, 0.004382848739624023, 0.9444444444444444)
    4  = @(0.0504612922668457, 0.01475977897644043, 0.9210526315789473)
    5  = @(0.1928677558898926, 0.09611129760742188, 0.7444444444444445)
    6  = @(0.1171343326568604, 0.03694820404052734, 0.9963636363636363)
    7  = @(3.9241783618927, 0.3028881549835205, 0.8889288281811646)
    8  = @(1.786609172821045, 0.06383275985717773, 0.9894827586206897)
    9  = @(2.114028215408325, 0.397730827331543, 0.8212713686150066)
    10 = @(19.88306164741516, 6.864698886871338, 0.877365)
    11 = @(0.1742157936096191, 0.04530668258666992, 7697.213483146067)
    12 = @(0.2116374969482422, 0.0632331371307373, 1783.098039215686)
    13 = @(0.2446317672729492, 0.1688055992126465, 0.9080902635702518)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 6).Value = $vals[0]
    $ws.Cells.Item($row, 7).Value = $vals[1]
    $ws.Cells.Item($row, 8).Value = $vals[2]
}
